$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.384.80'
$ws.Range("E2").Value = '  +2.74%  '
$ws.Range("D3").Value = '2.063.87'
$ws.Range("E3").Value = '  +4.39%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.65'
$ws.Range("E5").Value = '  +1.87%  '
$ws.Range("E6").Value = '  +3.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.07'
$ws.Range("E7").Value = '  +6.91%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  +3.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '57.87'
$ws.Range("E10").Value = '  -1.00%  '
$ws.Range("E11").Value = '  +1.92%  '
$ws.Range("E12").Value = '  +3.51%  '
$ws.Range("D13").Value = '2.367.94'
$ws.Range("E13").Value = '  +4.48%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.46'
$ws.Range("E14").Value = '  +4.07%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.00'
$ws.Range("E15").Value = '  +5.40%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.778'
$ws.Range("E16").Value = '  +3.89%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.20'
$ws.Range("E17").Value = '  +2.96%  '
$ws.Range("D18").Value = '2.061.16'
$ws.Range("E18").Value = '  +4.17%  '
$ws.Range("D19").Value = '37.589.15'
$ws.Range("E19").Value = '  +3.50%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.17'
$ws.Range("E20").Value = '  +16.90%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '69.22'
$ws.Range("E21").Value = '  +2.51%  '
$ws.Range("E22").Value = '  +1.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '227.93'
$ws.Range("E23").Value = '  +3.04%  '
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("E25").Value = '  +3.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.40'
$ws.Range("E26").Value = '  +1.71%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '163.79'
$ws.Range("E28").Value = '  +13.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.86'
$ws.Range("E29").Value = '  +3.80%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.21'
$ws.Range("E30").Value = '  +2.62%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.126'
$ws.Range("E31").Value = '  +1.77%  '
$ws.Range("E32").Value = '  +3.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.50'
$ws.Range("E33").Value = '  +3.51%  '
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.58'
$ws.Range("E34").Value = '  +12.34%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0622'
$ws.Range("E35").Value = '  +2.98%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.50'
$ws.Range("E36").Value = '  +6.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.43'
$ws.Range("E37").Value = '  +6.44%  '
$ws.Range("E38").Value = '  -0.17%  '
$ws.Range("E39").Value = '  +0.47%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.91'
$ws.Range("E40").Value = '  +11.61%  '
$ws.Range("B41").Value = 'Cronos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0991'
$ws.Range("E41").Value = '  +10.90%  '
$ws.Range("B42").Value = 'FTXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.58'
$ws.Range("E42").Value = '  +29.66%  '
$ws.Range("E43").Value = '  -1.81%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '97.23'
$ws.Range("E44").Value = '  +10.38%  '
$ws.Range("D45").Value = '1.476.54'
$ws.Range("E45").Value = '  +1.66%  '
$ws.Range("E46").Value = '  +7.75%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0211'
$ws.Range("E47").Value = '  +5.09%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.88'
$ws.Range("E48").Value = '  +7.60%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.03'
$ws.Range("E49").Value = '  +3.79%  '
$ws.Range("E50").Value = '  +6.71%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.94'
$ws.Range("E51").Value = '  +2.15%  '
